# CV: More wording improvements
# Applies the wording/content edits described by the commit to the active document.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replacement not found for:" $old
    }
}

function Get-ParagraphIndexContaining($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

# --- Simple whole-run text replacements -------------------------------------------------

Replace-Text "SPECIALTIES" "SKILLS"

Replace-Text "Managed releases: implemented a CI process" "Managed releases: established a CI process"

Replace-Text `
    "Played a key role in porting pilot software to iOS and Android, including successfully designing and implementing a mobile driving paradigm and UI that achieved both safety and usability." `
    "Played a key role in porting pilot software to iOS and Android, including successfully realizing a mobile driving model and UI that achieved both safety and usability."

Replace-Text `
    "Redesigned and rewrote the distributed scheduling architecture for device notifications." `
    "Redesigned and rewrote a distributed scheduling architecture for device notifications."

Replace-Text `
    "Internationalized code, automated text extraction/substitution, and added Korean support to the touch-screen keyboard (requiring me to learn the Hangul alphabet and all combining character logic)." `
    "Internationalized code, automated text extraction/substitution, and added Korean support to the touch-screen keyboard (requiring a crash course on the Hangul alphabet and combining character logic)."

Replace-Text `
    "Added HTML reporting, implemented as XSLT transforms on XML data." `
    "Provided HTML reporting via XSLT transforms on XML data."

Replace-Text `
    "Maintained a large part of the code for a BI reporting product, including UI design and implementation." `
    "Maintained much of the code for a BI reporting product."

Replace-Text `
    "Implemented a multi-dimensional OLAP reporting and interaction feature." `
    "Completed a multi-dimensional OLAP reporting and interaction feature."

Replace-Text `
    "Managed a team of developers: trained, delegated, instituted standards, and reviewed code." `
    "Directed a team of developers: trained, delegated, instituted standards, and reviewed code."

Replace-Text `
    "Architected the database and web UI for an automated asset auditing and management application." `
    "Architected the database and web UI for an asset auditing/management application."

Replace-Text `
    "homoiconic, embeddable programming language" `
    "homoiconic programming language"

Replace-Text `
    "and a single “program” data type." `
    "and a single “program” data type. A rewrite is in progress."

# --- Split the "neurofeedback" bullet into two bullets -----------------------------------
# Shorten the original bullet, then add a new sibling bullet (same list numbering) right
# after it holding the second half of the original sentence.

$neuroIdx = Get-ParagraphIndexContaining("Created a neurofeedback")
$neuroPara = $d.Paragraphs.Item($neuroIdx)
$neuroPara.Range.Text = "Created the BrainModder neurofeedback training software system."
$neuroPara.Range.InsertParagraphAfter()

$newBullet = $d.Paragraphs.Item($neuroIdx + 1)
$newBullet.Range.Text = "Developed networked games for BrainModder, controlled by the brain via real-time EEG and EMG data."

# --- Om Tree paragraph: several targeted in-run edits -------------------------------------

Replace-Text `
    "Creator of the Om Tree, an efficient associative array implemented in" `
    "Creator of the Om Tree, an efficient associative array written in"

Replace-Text `
    "was implemented for use in a forthcoming rewrite of the Om Language and has" `
    "has"

# Add a trailing empty run at the very end of the Om Tree paragraph (matching what a genuine
# Word paragraph split/merge leaves behind) by splitting right before the paragraph end and
# immediately re-joining it.
$omTreeIdx = Get-ParagraphIndexContaining("Creator of the Om Tree")
$omTreePara = $d.Paragraphs.Item($omTreeIdx)
$pEnd = $omTreePara.Range.End
$splitPoint = $d.Range($pEnd - 1, $pEnd - 1)
$splitPoint.InsertParagraphAfter()

$omTreePara2 = $d.Paragraphs.Item($omTreeIdx)
$mergeMark = $d.Range($omTreePara2.Range.End - 1, $omTreePara2.Range.End)
$mergeMark.Delete()

# Now rewrite the final "." run's text (without disturbing the hyperlink run before it).
$hyperlinkEnd = $d.Content
$hyperlinkEnd.Find.Execute("100% line, function, and branch unit test coverage") | Out-Null
$periodRun = $d.Range($hyperlinkEnd.End + 1, $hyperlinkEnd.End + 2)
$periodRun.Text = " and will be featured in the next version of the Om Language."
$periodRun.Font.Underline = 0
